$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.804.99'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '1.635.96'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''215.02'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = '''0.5071'
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '''0.2579'
$ws.Range("E8").Value = '  +0.60%  '
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("D10").Value = '''20.28'
$ws.Range("E10").Value = '  +4.04%  '
$ws.Range("D11").Value = '''0.07784'
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("D12").Value = '''4.249'
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '1.636.56'
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").Value = '1.862.78'
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = '''0.5595'
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Value = '0.0₅7648'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '''63.26'
$ws.Range("D18").Value = '25.803.01'
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("D19").Value = '''1.003'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").Value = '''4.373'
$ws.Range("E20").Value = '  -1.23%  '
$ws.Range("D21").Value = '''192.15'
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("D22").Value = '''9.901'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").Value = '''6.140'
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").Value = '''1.003'
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").Value = '''1.773'
$ws.Range("E25").Value = '  -6.21%  '
$ws.Range("D26").Value = '''138.94'
$ws.Range("E26").Value = '  -2.27%  '
$ws.Range("D27").Value = '''0.1230'
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("D28").Value = '''6.813'
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("D29").Value = '''15.52'
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").Value = '''1.242'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = '''0.04936'
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("D32").Value = '''3.287'
$ws.Range("E32").Value = '  +1.49%  '
$ws.Range("D33").Value = '''3.253'
$ws.Range("E33").Value = '  +2.06%  '
$ws.Range("D34").Value = '''1.571'
$ws.Range("E34").Value = '  +1.48%  '
$ws.Range("E35").Value = '  +0.45%  '
$ws.Range("D36").Value = '''0.9027'
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("E37").Value = '  +1.03%  '
$ws.Range("D38").Value = '''0.5560'
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("D39").Value = '1.131.54'
$ws.Range("E39").Value = '  +1.40%  '
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("D41").Value = '''0.9954'
$ws.Range("E41").Value = '  -0.59%  '
$ws.Range("D42").Value = '''5.448'
$ws.Range("E42").Value = '  -2.50%  '
$ws.Range("D43").Value = '''98.84'
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("D44").Value = '''0.7977'
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("D45").Value = '1.773.24'
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("E46").Value = '  -5.79%  '
$ws.Range("D47").Value = '''55.58'
$ws.Range("E47").Value = '  +1.31%  '
$ws.Range("D48").Value = '''0.4256'
$ws.Range("E48").Value = '  -4.09%  '
$ws.Range("D49").Value = '''7.745'
$ws.Range("E49").Value = '  +2.28%  '
$ws.Range("D50").Value = '''0.05030'
$ws.Range("D51").Value = '''0.9964'
$ws.Range("E51").Value = '  -0.71%  '
